$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1886.4849
$ws.Range("J40").Value = 1428.5
$ws.Range("L40").Value = 1428.5
$ws.Range("N40").Value = -1778.5
$ws.Range("H41").Value = 205
$ws.Range("I41").Value = 100
$ws.Range("J41").Value = 625
$ws.Range("K41").Value = 100
$ws.Range("L41").Value = 625
$ws.Range("M41").Value = 340
$ws.Range("N41").Value = -1505
$ws.Range("H64").Value = 367784.97
$ws.Range("I64").Value = 513368.94
$ws.Range("J64").Value = 3825
$ws.Range("K64").Value = 513368.94
$ws.Range("L64").Value = 3825
$ws.Range("M64").Value = -513120.94
$ws.Range("N64").Value = -4321
$ws.Range("H67").Value = 367784.97
$ws.Range("I67").Value = 513368.94
$ws.Range("J67").Value = 3825
$ws.Range("K67").Value = 513368.94
$ws.Range("L67").Value = 3825
$ws.Range("M67").Value = -512510.94
$ws.Range("N67").Value = -5541
$ws.Range("H74").Value = 5087
$ws.Range("I74").Value = 4376.143
$ws.Range("K74").Value = 4376.143
$ws.Range("M74").Value = -3440.143
$ws.Range("H76").Value = 3679.1777
$ws.Range("I76").Value = 3558.9429
$ws.Range("J76").Value = 4100
$ws.Range("K76").Value = 3558.9429
$ws.Range("L76").Value = 4100
$ws.Range("M76").Value = -3243.9429
$ws.Range("N76").Value = -4730
$ws.Range("H77").Value = 5087
$ws.Range("I77").Value = 4376.143
$ws.Range("K77").Value = 21880.715
$ws.Range("M77").Value = -17200.715
$ws.Range("H79").Value = 3679.1777
$ws.Range("I79").Value = 3558.9429
$ws.Range("J79").Value = 4100
$ws.Range("K79").Value = 3558.9429
$ws.Range("L79").Value = 4100
$ws.Range("M79").Value = -2466.9429
$ws.Range("N79").Value = -6284
$ws.Range("H112").Value = 6883.846
$ws.Range("I112").Value = 17815
$ws.Range("J112").Value = 2025.5555
$ws.Range("K112").Value = 53445
$ws.Range("L112").Value = 6076.666499999999
$ws.Range("M112").Value = -52337
$ws.Range("N112").Value = -8292.666499999999
$ws.Range("H116").Value = 2140.24
$ws.Range("I116").Value = 1920.6666
$ws.Range("J116").Value = 2469.6
$ws.Range("K116").Value = 1920.6666
$ws.Range("L116").Value = 2469.6
$ws.Range("M116").Value = 1521.3334
$ws.Range("N116").Value = -9353.6
$ws.Range("H128").Value = 39985.715
$ws.Range("J128").Value = 39985.715
$ws.Range("L128").Value = 39985.715
$ws.Range("N128").Value = -49945.715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3348.75
$ws.Range("I63").Value = 2798.3333
$ws.Range("K63").Value = 2798.3333
$ws.Range("M63").Value = -2112.3333
$ws.Range("H66").Value = 3348.75
$ws.Range("I66").Value = 2798.3333
$ws.Range("K66").Value = 13991.6665
$ws.Range("M66").Value = -10559.6665
$ws.Range("H114").Value = 30397
$ws.Range("J114").Value = 30397
$ws.Range("L114").Value = 30397
$ws.Range("N114").Value = -39075
$ws.Range("H122").Value = 12502244
$ws.Range("I122").Value = 2403.2856
$ws.Range("J122").Value = 41668536
$ws.Range("K122").Value = 7209.8568
$ws.Range("L122").Value = 125005608
$ws.Range("M122").Value = -4759.8568
$ws.Range("N122").Value = -125010508

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1973.6735
$ws.Range("I86").Value = 1802.3405
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 1802.3405
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -679.3405
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 1973.6735
$ws.Range("I89").Value = 1802.3405
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 9011.702499999999
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -3395.702499999999
$ws.Range("N89").Value = -41232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H117").Value = 62900
$ws.Range("J117").Value = 62900
$ws.Range("L117").Value = 62900
$ws.Range("N117").Value = -72078
$ws.Range("H123").Value = 75920
$ws.Range("J123").Value = 75920
$ws.Range("L123").Value = 75920
$ws.Range("N123").Value = -85720

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1094.8462
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 1166.6364
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 3499.9092
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -5121.9092
$ws.Range("H71").Value = 1094.8462
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 1166.6364
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 10499.7276
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -18611.7276
$ws.Range("H122").Value = 1073.8125
$ws.Range("I122").Value = 395.4
$ws.Range("J122").Value = 1382.1818
$ws.Range("K122").Value = 3558.6
$ws.Range("L122").Value = 12439.6362
$ws.Range("M122").Value = -1108.6
$ws.Range("N122").Value = -17339.6362
$ws.Range("H132").Value = 1248.8
$ws.Range("I132").Value = 931.3333
$ws.Range("J132").Value = 1725
$ws.Range("K132").Value = 8381.9997
$ws.Range("L132").Value = 15525
$ws.Range("M132").Value = -5851.9997
$ws.Range("N132").Value = -20585

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 25712.857
$ws.Range("J57").Value = 25712.857
$ws.Range("L57").Value = 25712.857
$ws.Range("N57").Value = -27352.857
$ws.Range("H70").Value = 5498.732
$ws.Range("I70").Value = 5367.304
$ws.Range("J70").Value = 5666.6665
$ws.Range("K70").Value = 5367.304
$ws.Range("L70").Value = 5666.6665
$ws.Range("M70").Value = -5097.304
$ws.Range("N70").Value = -6206.6665
$ws.Range("H73").Value = 5498.732
$ws.Range("I73").Value = 5367.304
$ws.Range("J73").Value = 5666.6665
$ws.Range("K73").Value = 5367.304
$ws.Range("L73").Value = 5666.6665
$ws.Range("M73").Value = -4431.304
$ws.Range("N73").Value = -7538.6665
$ws.Range("H103").Value = 300000
$ws.Range("J103").Value = 300000
$ws.Range("L103").Value = 300000
$ws.Range("N103").Value = -302344
$ws.Range("H122").Value = 7483.5386
$ws.Range("I122").Value = 11247.667
$ws.Range("J122").Value = 4257.143
$ws.Range("K122").Value = 33743.001
$ws.Range("L122").Value = 12771.429
$ws.Range("M122").Value = -31293.001
$ws.Range("N122").Value = -17671.429
$ws.Range("H128").Value = 57300
$ws.Range("J128").Value = 57300
$ws.Range("L128").Value = 57300
$ws.Range("N128").Value = -67260
$ws.Range("H129").Value = 39644.5
$ws.Range("J129").Value = 39644.5
$ws.Range("L129").Value = 39644.5
$ws.Range("N129").Value = -49644.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4941.0625
$ws.Range("I7").Value = 4677.909
$ws.Range("J7").Value = 5520
$ws.Range("K7").Value = 4677.909
$ws.Range("L7").Value = 5520
$ws.Range("M7").Value = -4565.909
$ws.Range("N7").Value = -5744
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = $null
$ws.Range("H40").Value = 2785.7693
$ws.Range("I40").Value = 2348.2
$ws.Range("J40").Value = 4244.3335
$ws.Range("K40").Value = 2348.2
$ws.Range("L40").Value = 4244.3335
$ws.Range("M40").Value = -2212.2
$ws.Range("N40").Value = -4516.3335
$ws.Range("H123").Value = 56533.332
$ws.Range("J123").Value = 56533.332
$ws.Range("L123").Value = 56533.332
$ws.Range("N123").Value = -66333.33199999999
$ws.Range("H126").Value = 4941.0625
$ws.Range("I126").Value = 4677.909
$ws.Range("J126").Value = 5520
$ws.Range("K126").Value = 14033.727
$ws.Range("L126").Value = 16560
$ws.Range("M126").Value = -11563.727
$ws.Range("N126").Value = -21500
$ws.Range("H132").Value = 4927.2856
$ws.Range("I132").Value = 4386.5557
$ws.Range("J132").Value = 5900.6
$ws.Range("K132").Value = 13159.6671
$ws.Range("L132").Value = 17701.8
$ws.Range("M132").Value = -10629.6671
$ws.Range("N132").Value = -22761.8
$ws.Range("H136").Value = 4548.7207
$ws.Range("I136").Value = 2199.739
$ws.Range("J136").Value = 7250.05
$ws.Range("K136").Value = 6599.217000000001
$ws.Range("L136").Value = 21750.15
$ws.Range("M136").Value = -4049.217000000001
$ws.Range("N136").Value = -26850.15

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6463.8
$ws.Range("I122").Value = 1662.3636
$ws.Range("J122").Value = 12332.223
$ws.Range("K122").Value = 4987.0908
$ws.Range("L122").Value = 36996.669
$ws.Range("M122").Value = -2537.0908
$ws.Range("N122").Value = -41896.669
$ws.Range("H132").Value = 3104.8635
$ws.Range("I132").Value = 3207
$ws.Range("J132").Value = 2832.5
$ws.Range("K132").Value = 9621
$ws.Range("L132").Value = 8497.5
$ws.Range("M132").Value = -7091
$ws.Range("N132").Value = -13557.5

Write-Host "Applied all profit sheet updates"